$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the course start date to C2 (store as literal text, not a date serial)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2024-08-29"
$ws.Range("C2").Style = "Normal"

# Add link to the ER random networks demo
$ws.Range("E9").Value = "[ER random networks demo](https://datahub.berkeley.edu/hub/user-redirect/git-pull?repo=https%3A%2F%2Fgithub.com%2Fdfeehan%2Fdemog180-fa2024&branch=main&urlpath=tree%2Fdemog180-fa2024%2Flecture%2F20240923_er_random_networks%2Fer_random_networks.ipynb)"

# Move the "Small worlds" section anchor up to row 10, and clear the old A11
$ws.Range("A10").Value = "[Small worlds](#sec:smallworlds)"
$ws.Range("A11").Value = ""

# Clear the leftover TODO demo placeholders
$ws.Range("E10").Value = ""
$ws.Range("E11").Value = ""

# Add link to Hwk 4 problem set
$ws.Range("G10").Value = "[Hwk 4: Problem set I](https://drive.google.com/file/d/1TUta8-8redraG0L044teOdA3SxX2eBtj/view?usp=sharing)"

# Add link to Lab 5 - Two-mode networks
$ws.Range("F11").Value = "[Lab 5 - Two-mode networks](https://datahub.berkeley.edu/hub/user-redirect/git-pull?repo=https%3A%2F%2Fgithub.com%2Fdfeehan%2Fdemog180-fa2024&branch=main&urlpath=tree%2Fdemog180-fa2024%2Flabs%2Flab05%2Flab05_two_mode_networks.ipynb)"
